# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Ajo" (garlic) at rows 282-283,
# shifting all subsequent records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 282 (pushes existing rows 282.. down by 2)
$ws.Rows.Item(282).EntireRow.Insert()
$ws.Rows.Item(282).EntireRow.Insert()

# New row 282
$ws.Cells.Item(282, 1).Value = 8
$ws.Cells.Item(282, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(282, 3).Value = "Coquimbo"
$ws.Cells.Item(282, 4).Value = 44855
$ws.Cells.Item(282, 5).Value = 4
$ws.Cells.Item(282, 6).Value = 100112003
$ws.Cells.Item(282, 7).Value = "Ajo"
$ws.Cells.Item(282, 8).Value = "Chino"
$ws.Cells.Item(282, 9).Value = "Primera"
$ws.Cells.Item(282, 10).Value = 520
$ws.Cells.Item(282, 11).Value = 16500
$ws.Cells.Item(282, 12).Value = 17000
$ws.Cells.Item(282, 13).Value = 16750
$ws.Cells.Item(282, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(282, 15).Value = "China"
$ws.Cells.Item(282, 16).Value = 1675
$ws.Cells.Item(282, 17).Value = 10
$ws.Cells.Item(282, 18).Value = "Hortaliza"

# New row 283
$ws.Cells.Item(283, 1).Value = 8
$ws.Cells.Item(283, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(283, 3).Value = "Coquimbo"
$ws.Cells.Item(283, 4).Value = 44855
$ws.Cells.Item(283, 5).Value = 4
$ws.Cells.Item(283, 6).Value = 100112003
$ws.Cells.Item(283, 7).Value = "Ajo"
$ws.Cells.Item(283, 8).Value = "Chino"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 500
$ws.Cells.Item(283, 11).Value = 19500
$ws.Cells.Item(283, 12).Value = 20000
$ws.Cells.Item(283, 13).Value = 19750
$ws.Cells.Item(283, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(283, 15).Value = "China"
$ws.Cells.Item(283, 16).Value = 1975
$ws.Cells.Item(283, 17).Value = 10
$ws.Cells.Item(283, 18).Value = "Hortaliza"
